# ================================================================
# Patient addition: add 4 new patients (100001-100004) representing
# tizanidine drug-drug-interaction (DDI) scenarios, each paired with
# a different interacting drug (Ciprofloxacin / zafirlukast /
# Fluvoxamine / Phenylpropanolamine). Mirrors the existing 6-row-per
# -patient block pattern (PER / VIS / D_ERA / D_EXP / D_EXP / D_ERA).
# ================================================================
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Step 1: write cell values for the template patient block (rows 701-706, patient 100001) ----
$ws.Range("A701").Value2 = "PER"
$ws.Range("B701").Value2 = 100001
$ws.Range("C701").Value2 = 1960
$ws.Range("K701").Value2 = 1
$ws.Range("L701").Value2 = 1
$ws.Range("O701").Value2 = 8532
$ws.Range("A702").Value2 = "VIS"
$ws.Range("B702").Value2 = 100001
$ws.Range("E702").Value2 = 44317
$ws.Range("F702").Value2 = 44323
$ws.Range("H702").Value2 = 9201
$ws.Range("A703").Value2 = "D_ERA"
$ws.Range("B703").Value2 = 100001
$ws.Range("D703").Value2 = 778474
$ws.Range("E703").Value2 = 44318
$ws.Range("F703").Value2 = 44318
$ws.Range("P703").Value2 = "Tizanidine"
$ws.Range("A704").Value2 = "D_EXP"
$ws.Range("B704").Value2 = 100001
$ws.Range("D704").Value2 = 778478
$ws.Range("E704").Value2 = 44319
$ws.Range("F704").Value2 = 44319
$ws.Range("I704").Value2 = 1
$ws.Range("J704").Value2 = 10
$ws.Range("M704").Value2 = "3 times daily"
$ws.Range("N704").Value2 = "null"
$ws.Range("P704").Value2 = "tizanidine 2 MG Oral Tablet"
$ws.Range("A705").Value2 = "D_EXP"
$ws.Range("B705").Value2 = 100001
$ws.Range("D705").Value2 = 19075391
$ws.Range("E705").Value2 = 44319
$ws.Range("F705").Value2 = 44319
$ws.Range("I705").Value2 = 1
$ws.Range("J705").Value2 = 10
$ws.Range("M705").Value2 = "3 times daily"
$ws.Range("N705").Value2 = "null"
$ws.Range("P705").Value2 = "Ciprofloxacin 100 MG/ML Oral Suspension"
$ws.Range("A706").Value2 = "D_ERA"
$ws.Range("B706").Value2 = 100001
$ws.Range("D706").Value2 = 1797513
$ws.Range("E706").Value2 = 44318
$ws.Range("F706").Value2 = 44318
$ws.Range("P706").Value2 = "Ciprofloxacin"

# ---- Step 2: apply correct cell styles to the template block (copied from existing cells with matching style) ----
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A701").PasteSpecial(-4122) | Out-Null
$ws.Range("B2").Copy() | Out-Null
$ws.Range("B701").PasteSpecial(-4122) | Out-Null
$ws.Range("C2").Copy() | Out-Null
$ws.Range("C701").PasteSpecial(-4122) | Out-Null
$ws.Range("D2").Copy() | Out-Null
$ws.Range("D701").PasteSpecial(-4122) | Out-Null
$ws.Range("E2").Copy() | Out-Null
$ws.Range("E701").PasteSpecial(-4122) | Out-Null
$ws.Range("F2").Copy() | Out-Null
$ws.Range("F701").PasteSpecial(-4122) | Out-Null
$ws.Range("G2").Copy() | Out-Null
$ws.Range("G701").PasteSpecial(-4122) | Out-Null
$ws.Range("H2").Copy() | Out-Null
$ws.Range("H701").PasteSpecial(-4122) | Out-Null
$ws.Range("I2").Copy() | Out-Null
$ws.Range("I701").PasteSpecial(-4122) | Out-Null
$ws.Range("J2").Copy() | Out-Null
$ws.Range("J701").PasteSpecial(-4122) | Out-Null
$ws.Range("K2").Copy() | Out-Null
$ws.Range("K701").PasteSpecial(-4122) | Out-Null
$ws.Range("L2").Copy() | Out-Null
$ws.Range("L701").PasteSpecial(-4122) | Out-Null
$ws.Range("M2").Copy() | Out-Null
$ws.Range("M701").PasteSpecial(-4122) | Out-Null
$ws.Range("N2").Copy() | Out-Null
$ws.Range("N701").PasteSpecial(-4122) | Out-Null
$ws.Range("O2").Copy() | Out-Null
$ws.Range("O701").PasteSpecial(-4122) | Out-Null
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A702").PasteSpecial(-4122) | Out-Null
$ws.Range("B2").Copy() | Out-Null
$ws.Range("B702").PasteSpecial(-4122) | Out-Null
$ws.Range("C3").Copy() | Out-Null
$ws.Range("C702").PasteSpecial(-4122) | Out-Null
$ws.Range("D2").Copy() | Out-Null
$ws.Range("D702").PasteSpecial(-4122) | Out-Null
$ws.Range("E3").Copy() | Out-Null
$ws.Range("E702").PasteSpecial(-4122) | Out-Null
$ws.Range("F3").Copy() | Out-Null
$ws.Range("F702").PasteSpecial(-4122) | Out-Null
$ws.Range("G2").Copy() | Out-Null
$ws.Range("G702").PasteSpecial(-4122) | Out-Null
$ws.Range("H6").Copy() | Out-Null
$ws.Range("H702").PasteSpecial(-4122) | Out-Null
$ws.Range("I2").Copy() | Out-Null
$ws.Range("I702").PasteSpecial(-4122) | Out-Null
$ws.Range("J2").Copy() | Out-Null
$ws.Range("J702").PasteSpecial(-4122) | Out-Null
$ws.Range("K3").Copy() | Out-Null
$ws.Range("K702").PasteSpecial(-4122) | Out-Null
$ws.Range("L3").Copy() | Out-Null
$ws.Range("L702").PasteSpecial(-4122) | Out-Null
$ws.Range("M2").Copy() | Out-Null
$ws.Range("M702").PasteSpecial(-4122) | Out-Null
$ws.Range("N2").Copy() | Out-Null
$ws.Range("N702").PasteSpecial(-4122) | Out-Null
$ws.Range("O3").Copy() | Out-Null
$ws.Range("O702").PasteSpecial(-4122) | Out-Null
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A703").PasteSpecial(-4122) | Out-Null
$ws.Range("B2").Copy() | Out-Null
$ws.Range("B703").PasteSpecial(-4122) | Out-Null
$ws.Range("C3").Copy() | Out-Null
$ws.Range("C703").PasteSpecial(-4122) | Out-Null
$ws.Range("D674").Copy() | Out-Null
$ws.Range("D703").PasteSpecial(-4122) | Out-Null
$ws.Range("E3").Copy() | Out-Null
$ws.Range("E703").PasteSpecial(-4122) | Out-Null
$ws.Range("F3").Copy() | Out-Null
$ws.Range("F703").PasteSpecial(-4122) | Out-Null
$ws.Range("G2").Copy() | Out-Null
$ws.Range("G703").PasteSpecial(-4122) | Out-Null
$ws.Range("H2").Copy() | Out-Null
$ws.Range("H703").PasteSpecial(-4122) | Out-Null
$ws.Range("I2").Copy() | Out-Null
$ws.Range("I703").PasteSpecial(-4122) | Out-Null
$ws.Range("J2").Copy() | Out-Null
$ws.Range("J703").PasteSpecial(-4122) | Out-Null
$ws.Range("K3").Copy() | Out-Null
$ws.Range("K703").PasteSpecial(-4122) | Out-Null
$ws.Range("L3").Copy() | Out-Null
$ws.Range("L703").PasteSpecial(-4122) | Out-Null
$ws.Range("M2").Copy() | Out-Null
$ws.Range("M703").PasteSpecial(-4122) | Out-Null
$ws.Range("N2").Copy() | Out-Null
$ws.Range("N703").PasteSpecial(-4122) | Out-Null
$ws.Range("O3").Copy() | Out-Null
$ws.Range("O703").PasteSpecial(-4122) | Out-Null
$ws.Range("P681").Copy() | Out-Null
$ws.Range("P703").PasteSpecial(-4122) | Out-Null
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A704").PasteSpecial(-4122) | Out-Null
$ws.Range("B2").Copy() | Out-Null
$ws.Range("B704").PasteSpecial(-4122) | Out-Null
$ws.Range("C3").Copy() | Out-Null
$ws.Range("C704").PasteSpecial(-4122) | Out-Null
$ws.Range("D674").Copy() | Out-Null
$ws.Range("D704").PasteSpecial(-4122) | Out-Null
$ws.Range("E3").Copy() | Out-Null
$ws.Range("E704").PasteSpecial(-4122) | Out-Null
$ws.Range("F3").Copy() | Out-Null
$ws.Range("F704").PasteSpecial(-4122) | Out-Null
$ws.Range("G2").Copy() | Out-Null
$ws.Range("G704").PasteSpecial(-4122) | Out-Null
$ws.Range("H2").Copy() | Out-Null
$ws.Range("H704").PasteSpecial(-4122) | Out-Null
$ws.Range("I7").Copy() | Out-Null
$ws.Range("I704").PasteSpecial(-4122) | Out-Null
$ws.Range("J7").Copy() | Out-Null
$ws.Range("J704").PasteSpecial(-4122) | Out-Null
$ws.Range("K3").Copy() | Out-Null
$ws.Range("K704").PasteSpecial(-4122) | Out-Null
$ws.Range("L3").Copy() | Out-Null
$ws.Range("L704").PasteSpecial(-4122) | Out-Null
$ws.Range("M7").Copy() | Out-Null
$ws.Range("M704").PasteSpecial(-4122) | Out-Null
$ws.Range("N7").Copy() | Out-Null
$ws.Range("N704").PasteSpecial(-4122) | Out-Null
$ws.Range("O3").Copy() | Out-Null
$ws.Range("O704").PasteSpecial(-4122) | Out-Null
$ws.Range("P683").Copy() | Out-Null
$ws.Range("P704").PasteSpecial(-4122) | Out-Null
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A705").PasteSpecial(-4122) | Out-Null
$ws.Range("B2").Copy() | Out-Null
$ws.Range("B705").PasteSpecial(-4122) | Out-Null
$ws.Range("C3").Copy() | Out-Null
$ws.Range("C705").PasteSpecial(-4122) | Out-Null
$ws.Range("D674").Copy() | Out-Null
$ws.Range("D705").PasteSpecial(-4122) | Out-Null
$ws.Range("E3").Copy() | Out-Null
$ws.Range("E705").PasteSpecial(-4122) | Out-Null
$ws.Range("F3").Copy() | Out-Null
$ws.Range("F705").PasteSpecial(-4122) | Out-Null
$ws.Range("G2").Copy() | Out-Null
$ws.Range("G705").PasteSpecial(-4122) | Out-Null
$ws.Range("H2").Copy() | Out-Null
$ws.Range("H705").PasteSpecial(-4122) | Out-Null
$ws.Range("I7").Copy() | Out-Null
$ws.Range("I705").PasteSpecial(-4122) | Out-Null
$ws.Range("J7").Copy() | Out-Null
$ws.Range("J705").PasteSpecial(-4122) | Out-Null
$ws.Range("K3").Copy() | Out-Null
$ws.Range("K705").PasteSpecial(-4122) | Out-Null
$ws.Range("L3").Copy() | Out-Null
$ws.Range("L705").PasteSpecial(-4122) | Out-Null
$ws.Range("M7").Copy() | Out-Null
$ws.Range("M705").PasteSpecial(-4122) | Out-Null
$ws.Range("N7").Copy() | Out-Null
$ws.Range("N705").PasteSpecial(-4122) | Out-Null
$ws.Range("O3").Copy() | Out-Null
$ws.Range("O705").PasteSpecial(-4122) | Out-Null
$ws.Range("P683").Copy() | Out-Null
$ws.Range("P705").PasteSpecial(-4122) | Out-Null
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A706").PasteSpecial(-4122) | Out-Null
$ws.Range("B2").Copy() | Out-Null
$ws.Range("B706").PasteSpecial(-4122) | Out-Null
$ws.Range("C3").Copy() | Out-Null
$ws.Range("C706").PasteSpecial(-4122) | Out-Null
$ws.Range("D674").Copy() | Out-Null
$ws.Range("D706").PasteSpecial(-4122) | Out-Null
$ws.Range("E3").Copy() | Out-Null
$ws.Range("E706").PasteSpecial(-4122) | Out-Null
$ws.Range("F3").Copy() | Out-Null
$ws.Range("F706").PasteSpecial(-4122) | Out-Null
$ws.Range("G2").Copy() | Out-Null
$ws.Range("G706").PasteSpecial(-4122) | Out-Null
$ws.Range("H2").Copy() | Out-Null
$ws.Range("H706").PasteSpecial(-4122) | Out-Null
$ws.Range("I2").Copy() | Out-Null
$ws.Range("I706").PasteSpecial(-4122) | Out-Null
$ws.Range("J2").Copy() | Out-Null
$ws.Range("J706").PasteSpecial(-4122) | Out-Null
$ws.Range("K3").Copy() | Out-Null
$ws.Range("K706").PasteSpecial(-4122) | Out-Null
$ws.Range("L3").Copy() | Out-Null
$ws.Range("L706").PasteSpecial(-4122) | Out-Null
$ws.Range("M2").Copy() | Out-Null
$ws.Range("M706").PasteSpecial(-4122) | Out-Null
$ws.Range("N2").Copy() | Out-Null
$ws.Range("N706").PasteSpecial(-4122) | Out-Null
$ws.Range("O3").Copy() | Out-Null
$ws.Range("O706").PasteSpecial(-4122) | Out-Null
$ws.Range("P683").Copy() | Out-Null
$ws.Range("P706").PasteSpecial(-4122) | Out-Null

# ---- Step 3: replicate the formatted template block for the other 3 patients ----
$ws.Range("A701:O706").Copy() | Out-Null
$ws.Range("A707:O712").PasteSpecial(-4122) | Out-Null
$ws.Range("A713:O718").PasteSpecial(-4122) | Out-Null
$ws.Range("A719:O724").PasteSpecial(-4122) | Out-Null

# ---- Step 4: write cell values that differ per patient (PERSON_ID + per-patient drug cells) ----
$ws.Range("A707").Value2 = "PER"
$ws.Range("B707").Value2 = 100002
$ws.Range("C707").Value2 = 1960
$ws.Range("K707").Value2 = 1
$ws.Range("L707").Value2 = 1
$ws.Range("O707").Value2 = 8532
$ws.Range("A708").Value2 = "VIS"
$ws.Range("B708").Value2 = 100002
$ws.Range("E708").Value2 = 44317
$ws.Range("F708").Value2 = 44323
$ws.Range("H708").Value2 = 9201
$ws.Range("A709").Value2 = "D_ERA"
$ws.Range("B709").Value2 = 100002
$ws.Range("D709").Value2 = 778474
$ws.Range("E709").Value2 = 44318
$ws.Range("F709").Value2 = 44318
$ws.Range("P709").Value2 = "Tizanidine"
$ws.Range("A710").Value2 = "D_EXP"
$ws.Range("B710").Value2 = 100002
$ws.Range("D710").Value2 = 778478
$ws.Range("E710").Value2 = 44319
$ws.Range("F710").Value2 = 44319
$ws.Range("I710").Value2 = 1
$ws.Range("J710").Value2 = 10
$ws.Range("M710").Value2 = "3 times daily"
$ws.Range("N710").Value2 = "null"
$ws.Range("P710").Value2 = "tizanidine 2 MG Oral Tablet"
$ws.Range("A711").Value2 = "D_EXP"
$ws.Range("B711").Value2 = 100002
$ws.Range("D711").Value2 = 1111710
$ws.Range("E711").Value2 = 44319
$ws.Range("F711").Value2 = 44319
$ws.Range("I711").Value2 = 1
$ws.Range("J711").Value2 = 10
$ws.Range("M711").Value2 = "3 times daily"
$ws.Range("N711").Value2 = "null"
$ws.Range("P711").Value2 = "zafirlukast 10 MG Oral Tablet"
$ws.Range("A712").Value2 = "D_ERA"
$ws.Range("B712").Value2 = 100002
$ws.Range("D712").Value2 = 1111706
$ws.Range("E712").Value2 = 44318
$ws.Range("F712").Value2 = 44318
$ws.Range("P712").Value2 = "zafirlukast"
$ws.Range("A713").Value2 = "PER"
$ws.Range("B713").Value2 = 100003
$ws.Range("C713").Value2 = 1960
$ws.Range("K713").Value2 = 1
$ws.Range("L713").Value2 = 1
$ws.Range("O713").Value2 = 8532
$ws.Range("A714").Value2 = "VIS"
$ws.Range("B714").Value2 = 100003
$ws.Range("E714").Value2 = 44317
$ws.Range("F714").Value2 = 44323
$ws.Range("H714").Value2 = 9201
$ws.Range("A715").Value2 = "D_ERA"
$ws.Range("B715").Value2 = 100003
$ws.Range("D715").Value2 = 778474
$ws.Range("E715").Value2 = 44318
$ws.Range("F715").Value2 = 44318
$ws.Range("P715").Value2 = "Tizanidine"
$ws.Range("A716").Value2 = "D_EXP"
$ws.Range("B716").Value2 = 100003
$ws.Range("D716").Value2 = 778478
$ws.Range("E716").Value2 = 44319
$ws.Range("F716").Value2 = 44319
$ws.Range("I716").Value2 = 1
$ws.Range("J716").Value2 = 10
$ws.Range("M716").Value2 = "3 times daily"
$ws.Range("N716").Value2 = "null"
$ws.Range("P716").Value2 = "tizanidine 2 MG Oral Tablet"
$ws.Range("A717").Value2 = "D_EXP"
$ws.Range("B717").Value2 = 100003
$ws.Range("D717").Value2 = 40174735
$ws.Range("E717").Value2 = 44319
$ws.Range("F717").Value2 = 44319
$ws.Range("I717").Value2 = 1
$ws.Range("J717").Value2 = 10
$ws.Range("M717").Value2 = "3 times daily"
$ws.Range("N717").Value2 = "null"
$ws.Range("P717").Value2 = "Fluvoxamine Maleate 25 MG Oral Tablet"
$ws.Range("A718").Value2 = "D_ERA"
$ws.Range("B718").Value2 = 100003
$ws.Range("D718").Value2 = 751412
$ws.Range("E718").Value2 = 44318
$ws.Range("F718").Value2 = 44318
$ws.Range("P718").Value2 = "Fluvoxamine"
$ws.Range("A719").Value2 = "PER"
$ws.Range("B719").Value2 = 100004
$ws.Range("C719").Value2 = 1960
$ws.Range("K719").Value2 = 1
$ws.Range("L719").Value2 = 1
$ws.Range("O719").Value2 = 8532
$ws.Range("A720").Value2 = "VIS"
$ws.Range("B720").Value2 = 100004
$ws.Range("E720").Value2 = 44317
$ws.Range("F720").Value2 = 44323
$ws.Range("H720").Value2 = 9201
$ws.Range("A721").Value2 = "D_ERA"
$ws.Range("B721").Value2 = 100004
$ws.Range("D721").Value2 = 778474
$ws.Range("E721").Value2 = 44318
$ws.Range("F721").Value2 = 44318
$ws.Range("P721").Value2 = "Tizanidine"
$ws.Range("A722").Value2 = "D_EXP"
$ws.Range("B722").Value2 = 100004
$ws.Range("D722").Value2 = 778478
$ws.Range("E722").Value2 = 44319
$ws.Range("F722").Value2 = 44319
$ws.Range("I722").Value2 = 1
$ws.Range("J722").Value2 = 10
$ws.Range("M722").Value2 = "3 times daily"
$ws.Range("N722").Value2 = "null"
$ws.Range("P722").Value2 = "tizanidine 2 MG Oral Tablet"
$ws.Range("A723").Value2 = "D_EXP"
$ws.Range("B723").Value2 = 100004
$ws.Range("D723").Value2 = 40243565
$ws.Range("E723").Value2 = 44319
$ws.Range("F723").Value2 = 44319
$ws.Range("I723").Value2 = 1
$ws.Range("J723").Value2 = 10
$ws.Range("M723").Value2 = "3 times daily"
$ws.Range("N723").Value2 = "null"
$ws.Range("P723").Value2 = "Phenylpropanolamine Hydrochloride 50 MG Chewable Tablet"
$ws.Range("A724").Value2 = "D_ERA"
$ws.Range("B724").Value2 = 100004
$ws.Range("D724").Value2 = 1139993
$ws.Range("E724").Value2 = 44318
$ws.Range("F724").Value2 = 44318
$ws.Range("P724").Value2 = "Phenylpropanolamine"

# ---- Step 5: set the active selection to match the authored workbook view state ----
$ws.Range("F713").Select() | Out-Null
